$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Edn3"
$ws.Range("C2").Value = "Ednra"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.012975
$ws.Range("H2").Value = 0.038925
$ws.Range("I2").Value = 0.004878150260562778
$ws.Range("J2").Value = 0.004878150260562778
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 6.433013333333332
$ws.Range("N2").Value = 19.29904
$ws.Range("O2").Value = 0.1097146002786867
$ws.Range("P2").Value = 0.1097146002786867
$ws.Range("Q2").Value = 0.08346834799999998
$ws.Range("R2").Value = 0.7512151319999999
$ws.Range("S2").Value = 0.0005352043059370163
$ws.Range("T2").Value = 0.0005352043059370165

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Edn3"
$ws.Range("C3").Value = "Ednra"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.012975
$ws.Range("H3").Value = 0.038925
$ws.Range("I3").Value = 0.004878150260562778
$ws.Range("J3").Value = 0.004878150260562778
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 32.709374
$ws.Range("N3").Value = 98.12812199999999
$ws.Range("O3").Value = 0.5578561255548565
$ws.Range("P3").Value = 0.5578561255548566
$ws.Range("Q3").Value = 0.42440412765
$ws.Range("R3").Value = 3.81963714885
$ws.Range("S3").Value = 0.002721306004231965
$ws.Range("T3").Value = 0.002721306004231966

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Edn3"
$ws.Range("C4").Value = "Ednra"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.012975
$ws.Range("H4").Value = 0.038925
$ws.Range("I4").Value = 0.004878150260562778
$ws.Range("J4").Value = 0.004878150260562778
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.07263
$ws.Range("N4").Value = 0.21789
$ws.Range("O4").Value = 0.001238699658362439
$ws.Range("P4").Value = 0.001238699658362439
$ws.Range("Q4").Value = 0.00094237425
$ws.Range("R4").Value = 0.008481368250000001
$ws.Range("S4").Value = 0.000006042563061199755
$ws.Range("T4").Value = 0.000006042563061199755

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Edn3"
$ws.Range("C5").Value = "Ednra"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.012975
$ws.Range("H5").Value = 0.038925
$ws.Range("I5").Value = 0.004878150260562778
$ws.Range("J5").Value = 0.004878150260562778
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 19.41905066666667
$ws.Range("N5").Value = 58.257152
$ws.Range("O5").Value = 0.3311905745080943
$ws.Range("P5").Value = 0.3311905745080943
$ws.Range("Q5").Value = 0.2519621824
$ws.Range("R5").Value = 2.2676596416
$ws.Range("S5").Value = 0.001615597387332596
$ws.Range("T5").Value = 0.001615597387332596

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Edn3"
$ws.Range("C6").Value = "Ednra"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.646844666666667
$ws.Range("H6").Value = 7.940534
$ws.Range("I6").Value = 0.9951218497394373
$ws.Range("J6").Value = 0.9951218497394373
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 6.433013333333332
$ws.Range("N6").Value = 19.29904
$ws.Range("O6").Value = 0.1097146002786867
$ws.Range("P6").Value = 0.1097146002786867
$ws.Range("Q6").Value = 17.02718703192889
$ws.Range("R6").Value = 153.24468328736
$ws.Range("S6").Value = 0.1091793959727496
$ws.Range("T6").Value = 0.1091793959727497

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Edn3"
$ws.Range("C7").Value = "Ednra"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.646844666666667
$ws.Range("H7").Value = 7.940534
$ws.Range("I7").Value = 0.9951218497394373
$ws.Range("J7").Value = 0.9951218497394373
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 32.709374
$ws.Range("N7").Value = 98.12812199999999
$ws.Range("O7").Value = 0.5578561255548565
$ws.Range("P7").Value = 0.5578561255548566
$ws.Range("Q7").Value = 86.57663212190532
$ws.Range("R7").Value = 779.1896890971478
$ws.Range("S7").Value = 0.5551348195506246
$ws.Range("T7").Value = 0.5551348195506247

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Edn3"
$ws.Range("C8").Value = "Ednra"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2.646844666666667
$ws.Range("H8").Value = 7.940534
$ws.Range("I8").Value = 0.9951218497394373
$ws.Range("J8").Value = 0.9951218497394373
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.07263
$ws.Range("N8").Value = 0.21789
$ws.Range("O8").Value = 0.001238699658362439
$ws.Range("P8").Value = 0.001238699658362439
$ws.Range("Q8").Value = 0.19224032814
$ws.Range("R8").Value = 1.73016295326
$ws.Range("S8").Value = 0.001232657095301239
$ws.Range("T8").Value = 0.001232657095301239

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Edn3"
$ws.Range("C9").Value = "Ednra"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2.646844666666667
$ws.Range("H9").Value = 7.940534
$ws.Range("I9").Value = 0.9951218497394373
$ws.Range("J9").Value = 0.9951218497394373
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 19.41905066666667
$ws.Range("N9").Value = 58.257152
$ws.Range("O9").Value = 0.3311905745080943
$ws.Range("P9").Value = 0.3311905745080943
$ws.Range("Q9").Value = 51.39921068879644
$ws.Range("R9").Value = 462.592896199168
$ws.Range("S9").Value = 0.3295749771207617
$ws.Range("T9").Value = 0.3295749771207617

